$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.841.25"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.902.91"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5043"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3814"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07293"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9097"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07656"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "1.873.77"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.478"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.604"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008694"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "27.874.10"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.156"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.858"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.235"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.28%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.920"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08991"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.216"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.239"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7693"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.637"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5543"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.015"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05260"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.515"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1526"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4795"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.636"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06080"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9012"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
